$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AD2").Value = 1.331617765164322
$ws.Range("AD3").Value = 0.1615255895249083
$ws.Range("AD4").Value = -0.2190164776938517
$ws.Range("AD5").Value = 1.247005318008205
$ws.Range("AD7").Value = 0.6325998303579028
$ws.Range("AD8").Value = 0.01465343221456516
$ws.Range("AD9").Value = 0.2522790160453459
$ws.Range("AD10").Value = 0.3692449490957154
$ws.Range("AD11").Value = 0.2401404958660949
$ws.Range("AD12").Value = 0.04541607648193594
$ws.Range("AD13").Value = 0.3439855200790409
$ws.Range("AD14").Value = -0.2771027461821775
$ws.Range("AD15").Value = 0.9423268497684154
$ws.Range("AD17").Value = 0.4934238596771008
$ws.Range("AD18").Value = 0.1819533668932887
$ws.Range("AD19").Value = 0.1274835671478239
$ws.Range("AD21").Value = 0.1175831222413791
$ws.Range("AD22").Value = 0.08567806172022781
$ws.Range("AD23").Value = 0.2472168566298102
$ws.Range("AD24").Value = 0.9657324797925104
$ws.Range("AD25").Value = 0.6240851513992403
$ws.Range("AD26").Value = 0.6257639271635584
$ws.Range("AD27").Value = 0.7119843018245802
$ws.Range("AD28").Value = 1.214379302193214
$ws.Range("AD29").Value = -0.2762158054446322
$ws.Range("AD30").Value = -1.416062347784709
$ws.Range("AD31").Value = 0.6333022396118085
$ws.Range("AD32").Value = -1.029959203170345
$ws.Range("AD33").Value = 0.6096467109995476
$ws.Range("AD34").Value = 0.160715332494081
$ws.Range("AD35").Value = 0.7917059169459646
$ws.Range("AD37").Value = -0.01968783795405432
$ws.Range("AD42").Value = 0.492203158203061
$ws.Range("AD43").Value = 1.166045041745317
$ws.Range("AD44").Value = 0.03885572782721866
$ws.Range("AD46").Value = 0.6781037513389003
$ws.Range("AD47").Value = 0.3439146509921893
$ws.Range("AD48").Value = -0.00429677831595654
$ws.Range("AD50").Value = 0.6832196971504075
$ws.Range("AD51").Value = -0.03057912132167632
$ws.Range("AD52").Value = 0.4212951720538849
$ws.Range("AD53").Value = -0.1357650071055301
$ws.Range("AD54").Value = -0.0347585287505181
$ws.Range("AD55").Value = 0.05214779532647944
$ws.Range("AD56").Value = -0.274777650061912
$ws.Range("AD57").Value = 0.1493591001217888
$ws.Range("AD58").Value = -0.5265188872074642
$ws.Range("AD59").Value = -0.6537420480262576
$ws.Range("AD60").Value = 0.1895020719705186
